# "Fix up Yolanda and Inst. Curie"
#
# The collaboration list currently has 4 data rows (rows 2-5):
#   row 2 - Yolanda Prezado (Dr.)
#   row 3 - Alfredo Fernandez-Rodriguez (M)
#   row 4 - Thongchai Masilela (Dr.)
#   row 5 - Frederic Pouzoulet (Dr.)
#
# The fix removes the Yolanda Prezado, Thongchai Masilela and Frederic
# Pouzoulet rows, leaving only the Alfredo Fernandez-Rodriguez row (which
# ends up as the sole data row, row 2) in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Yolanda Prezado row. This shifts every row below it up by
# one, so Alfredo's former row 3 becomes row 2.
$ws.Rows("2").Delete()

# Remove the Thongchai Masilela and Frederic Pouzoulet rows, which - after
# the shift above - are now rows 3 and 4.
$ws.Rows("3:4").Delete()

# Leave the selection where the author last left it.
$ws.Range("C10").Select()
